# Update "想去人数" (want-to-go count) values in the 展览 (Exhibition) sheet
# and the mirrored 全部类型 (All Types) sheet, reflecting newly generated
# output data (gh-pages build at 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# row -> new F-column value
$updates = @{
    2  = 259
    4  = 13602
    9  = 172
    10 = 248
    14 = 66
    16 = 49
    18 = 5635
    20 = 73
    22 = 51
    25 = 184
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}

$wb.Save()
